$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (preserve rich-text runs via Characters) ---

# A8: "Volume 32   Number  48" -> "...49"  (run "48" is chars 21-22)
$volCell = $ws.Range("A8")
$volCell.Characters(21, 2).Text = "49"

# C9: "Report Covering the Week  11/24/2025  Through  11/30/2025"
#     -> "...  12/1/2025  Through  12/7/2025"
$weekCell = $ws.Range("C9")
$weekCell.Characters(27, 10).Text = "12/1/2025"
$weekCell.Characters(47, 10).Text = "12/7/2025"

# --- Crime-statistics table updates (rows 15-31) ---

$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 7
$ws.Range("H15").Value = -57.142857142857
$ws.Range("I15").Value = 38
$ws.Range("J15").Value = 36
$ws.Range("K15").Value = 5.555555555555
$ws.Range("L15").Value = 15.151515151515
$ws.Range("M15").Value = 137.5
$ws.Range("N15").Value = 31.03448275862
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 21
$ws.Range("H16").Value = -28.571428571428
$ws.Range("I16").Value = 230
$ws.Range("J16").Value = 281
$ws.Range("K16").Value = -18.14946619217
$ws.Range("L16").Value = -25.324675324675
$ws.Range("M16").Value = 17.34693877551
$ws.Range("N16").Value = -79.717813051146
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 11
$ws.Range("E17").Value = -54.545454545454
$ws.Range("F17").Value = 25
$ws.Range("G17").Value = 36
$ws.Range("H17").Value = -30.555555555555
$ws.Range("I17").Value = 363
$ws.Range("J17").Value = 472
$ws.Range("K17").Value = -23.093220338983
$ws.Range("L17").Value = -13.157894736842
$ws.Range("M17").Value = 75.362318840579
$ws.Range("N17").Value = 11.349693251533
$ws.Range("C18").Value = 8
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 166.666666666667
$ws.Range("F18").Value = 25
$ws.Range("G18").Value = 36
$ws.Range("H18").Value = -30.555555555555
$ws.Range("I18").Value = 422
$ws.Range("J18").Value = 403
$ws.Range("K18").Value = 4.714640198511
$ws.Range("L18").Value = -19.77186311787
$ws.Range("M18").Value = -6.843267108167
$ws.Range("N18").Value = -81.244444444444
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 21
$ws.Range("E19").Value = -23.809523809523
$ws.Range("F19").Value = 75
$ws.Range("G19").Value = 84
$ws.Range("H19").Value = -10.714285714285
$ws.Range("I19").Value = 978
$ws.Range("J19").Value = 1214
$ws.Range("K19").Value = -19.439868204283
$ws.Range("L19").Value = -23.354231974921
$ws.Range("M19").Value = 53.77358490566
$ws.Range("N19").Value = -24.068322981366
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 40
$ws.Range("F20").Value = 41
$ws.Range("G20").Value = 30
$ws.Range("H20").Value = 36.666666666666
$ws.Range("I20").Value = 552
$ws.Range("J20").Value = 499
$ws.Range("K20").Value = 10.621242484969
$ws.Range("L20").Value = 14.522821576763
$ws.Range("M20").Value = 97.849462365591
$ws.Range("N20").Value = -85.809768637532
$ws.Range("C21").Value = 40
$ws.Range("D21").Value = 43
$ws.Range("E21").Value = -6.976744186046
$ws.Range("F21").Value = 184
$ws.Range("G21").Value = 214
$ws.Range("H21").Value = -14.018691588785
$ws.Range("I21").Value = 2586
$ws.Range("J21").Value = 2910
$ws.Range("K21").Value = -11.134020618556
$ws.Range("L21").Value = -15.101772816808
$ws.Range("M21").Value = 44.147157190635
$ws.Range("N21").Value = -71.034946236559
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -50
$ws.Range("I22").Value = 31
$ws.Range("J22").Value = 30
$ws.Range("K22").Value = 3.333333333333
$ws.Range("L22").Value = -6.060606060606
$ws.Range("M22").Value = 675
$ws.Range("C24").Value = 82
$ws.Range("D24").Value = 64
$ws.Range("E24").Value = 28.125
$ws.Range("F24").Value = 222
$ws.Range("G24").Value = 227
$ws.Range("H24").Value = -2.202643171806
$ws.Range("I24").Value = 2445
$ws.Range("J24").Value = 2614
$ws.Range("K24").Value = -6.46518745218
$ws.Range("L24").Value = -7.874905802562
$ws.Range("M24").Value = 64.868509777478
$ws.Range("C25").Value = 52
$ws.Range("D25").Value = 36
$ws.Range("E25").Value = 44.444444444444
$ws.Range("F25").Value = 113
$ws.Range("G25").Value = 138
$ws.Range("H25").Value = -18.115942028985
$ws.Range("I25").Value = 1393
$ws.Range("J25").Value = 1709
$ws.Range("K25").Value = -18.490345231129
$ws.Range("L25").Value = -5.623306233062
$ws.Range("C26").Value = 16
$ws.Range("D26").Value = 18
$ws.Range("E26").Value = -11.111111111111
$ws.Range("F26").Value = 64
$ws.Range("G26").Value = 83
$ws.Range("H26").Value = -22.89156626506
$ws.Range("I26").Value = 815
$ws.Range("J26").Value = 921
$ws.Range("K26").Value = -11.509229098805
$ws.Range("L26").Value = -2.628434886499
$ws.Range("M26").Value = 20.562130177514
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -28.571428571428
$ws.Range("I27").Value = 45
$ws.Range("J27").Value = 44
$ws.Range("K27").Value = 2.272727272727
$ws.Range("L27").Value = -10
$ws.Range("C28").Value = 3
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = 16.666666666666
$ws.Range("I28").Value = 77
$ws.Range("J28").Value = 94
$ws.Range("K28").Value = -18.085106382978
$ws.Range("L28").Value = -18.085106382978
$ws.Range("F29").Value = 1
$ws.Range("N29").Value = -69.230769230769
$ws.Range("F30").Value = 1
$ws.Range("N30").Value = -63.636363636363
$ws.Range("J31").Value = 11
$ws.Range("K31").Value = 9.090909090909

# --- Cells whose data type / number format changes (text <-> number) ---

$c = $ws.Range("C22")
$c.NumberFormat = '#,##0'
$c.Value = 2

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0"
$c.NumberFormat = "General"

$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "***.*"
$c.NumberFormat = "General"

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0"
$c.NumberFormat = "General"

$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "***.*"
$c.NumberFormat = "General"

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0"
$c.NumberFormat = "General"

$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = "***.*"
$c.NumberFormat = "General"

$c = $ws.Range("D31")
$c.NumberFormat = '#,##0'
$c.Value = 1

$c = $ws.Range("E31")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = -100

$c = $ws.Range("G31")
$c.NumberFormat = '#,##0'
$c.Value = 1

$c = $ws.Range("H31")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = 0

